$wb = $excel.ActiveWorkbook

# Delete the "Data Texas" worksheet entirely.
[void]$wb.Worksheets.Item("Data Texas").Delete()

# Update the base excess-capacity value on the HPPECbP sheet from 0.1 back to 0.25.
# (All other cells on the sheet reference this cell via formulas, so they update too.)
$ws = $wb.Worksheets.Item("HPPECbP")
$ws.Range("B2").Value = 0.25
